$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title: consolidate "Testing"/" "/"custom"/" "/"properties" runs into one run.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleText = $titleRange.Text
$titleRange.Characters(1, $titleText.Length).Text = $titleText

# --- Subtitle: consolidate the run-per-word text segments that sit between
#     the <a:br/> line breaks (code 11 / vertical tab), without touching the
#     breaks themselves.
$subRange = $s.Shapes.Item(2).TextFrame.TextRange
$subText = $subRange.Text

$segStart = 1
for ($i = 1; $i -le $subText.Length; $i++) {
    $code = [int][char]$subText[$i - 1]
    if ($code -eq 11) {
        if ($i -gt $segStart) {
            $segText = $subText.Substring($segStart - 1, $i - $segStart)
            $subRange.Characters($segStart, $i - $segStart).Text = $segText
        }
        $segStart = $i + 1
    }
}
if ($subText.Length -ge $segStart) {
    $segText = $subText.Substring($segStart - 1, $subText.Length - $segStart + 1)
    $subRange.Characters($segStart, $subText.Length - $segStart + 1).Text = $segText
}
